# Applies the dac-val.xlsx testdata edit:
#  - Replace the "Param2"/"Param3" column headers (D1/E1) with new,
#    arbitrarily-named parameter headers ("foo"/"bar"), reflecting the
#    move to supporting any number of parameters with any name.
#  - Update the per-row parameter values in columns D and E to the new
#    (normalized, 0-1 range) breakpoint values.
#  - Move the active selection to F11 (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parameter column headers (row 1)
$ws.Range("D1").Value = "foo"
$ws.Range("E1").Value = "bar"

# Row 2 (DSApplause)
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = 0.1

# Row 3 (DSBugs)
$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 0.9

# Row 4 (DSWind)
$ws.Range("D4").Value = 0.9
$ws.Range("E4").Value = 0.9

# Row 5 (DSPistons)
$ws.Range("D5").Value = 0.9
$ws.Range("E5").Value = 0.2

# Move the selection/active cell as recorded in the saved workbook view
$ws.Range("F11").Select()
